# Scheduled market-data refresh: update cached currentAveragePrice* / LevePrice* /
# LeveProfit* figures (columns H-N) for the affected leve rows on each job sheet.
# Source values come from a recalculated market snapshot; only these specific cells
# change (other leve rows / sheet layout are left untouched).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 487.66666
$ws.Range("I98").Value = 487.66666
$ws.Range("K98").Value = 487.66666
$ws.Range("M98").Value = 1010.33334

# Row 106: Making Your Mark / Enchanted Palladium Ink
$ws.Range("H106").Value = 1698.2667
$ws.Range("I106").Value = 1424
$ws.Range("J106").Value = 2452.5
$ws.Range("K106").Value = 1424
$ws.Range("L106").Value = 2452.5
$ws.Range("M106").Value = -793
$ws.Range("N106").Value = -3714.5

# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 1098.5862
$ws.Range("J112").Value = 1119.9642
$ws.Range("L112").Value = 3359.8926
$ws.Range("N112").Value = -5575.892599999999

# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 487.66666
$ws.Range("I122").Value = 487.66666
$ws.Range("K122").Value = 1462.99998
$ws.Range("M122").Value = 987.0000199999999

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 5356.4
$ws.Range("I132").Value = 6349.091
$ws.Range("K132").Value = 19047.273
$ws.Range("M132").Value = -16517.273

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1630.3518
$ws.Range("J138").Value = 2437.7334
$ws.Range("L138").Value = 7313.2002
$ws.Range("N138").Value = -17593.2002

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 18138.254
$ws.Range("I32").Value = 18726.203
$ws.Range("J32").Value = 9466
$ws.Range("K32").Value = 18726.203
$ws.Range("L32").Value = 9466
$ws.Range("M32").Value = -18439.203
$ws.Range("N32").Value = -10040

# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 3862.4119
$ws.Range("I45").Value = 3674
$ws.Range("K45").Value = 3674
$ws.Range("M45").Value = -3297

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 1663.625
$ws.Range("J102").Value = 1999.6666
$ws.Range("L102").Value = 1999.6666
$ws.Range("N102").Value = -5243.6666

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 14028.125
$ws.Range("I132").Value = 1228.0938
$ws.Range("K132").Value = 3684.2814
$ws.Range("M132").Value = -1154.2814

$ws = $wb.Worksheets.Item("BSM")
# Row 19: Twice as Slice / Spiked Bronze Labrys
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = ""

# Row 94: High Steal / High Steel Nugget
$ws.Range("H94").Value = 2261.3948
$ws.Range("I94").Value = 1064.4333
$ws.Range("K94").Value = 1064.4333
$ws.Range("M94").Value = -613.4332999999999

# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 1861.6111
$ws.Range("I99").Value = 1522.1111
$ws.Range("J99").Value = 2201.111
$ws.Range("K99").Value = 1522.1111
$ws.Range("L99").Value = 2201.111
$ws.Range("M99").Value = -24.11110000000008
$ws.Range("N99").Value = -5197.111

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 28296.3
$ws.Range("I134").Value = 41328.074
$ws.Range("K134").Value = 123984.222
$ws.Range("M134").Value = -121449.222

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 12623.793
$ws.Range("I31").Value = 15808.682
$ws.Range("K31").Value = 15808.682
$ws.Range("M31").Value = -15513.682

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 12623.793
$ws.Range("I34").Value = 15808.682
$ws.Range("K34").Value = 15808.682
$ws.Range("M34").Value = -15606.682

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 26289.045
$ws.Range("I132").Value = 31843.941
$ws.Range("K132").Value = 95531.823
$ws.Range("M132").Value = -93001.823

$ws = $wb.Worksheets.Item("CUL")
# Row 11: Putting the Squeeze On / Orange Juice
$ws.Range("H11").Value = 118.125
$ws.Range("I11").Value = 57.5
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 172.5
$ws.Range("L11").Value = 900
$ws.Range("M11").Value = -32.5
$ws.Range("N11").Value = -1180

# Row 124: Bobbing for Compliments / Island Miq'abob
$ws.Range("H124").Value = 1000
$ws.Range("I124").Value = 1000
$ws.Range("K124").Value = 3000
$ws.Range("M124").Value = 1910

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 115716.71
$ws.Range("I131").Value = 734.5
$ws.Range("J131").Value = 124233.914
$ws.Range("K131").Value = 2203.5
$ws.Range("L131").Value = 372701.742
$ws.Range("M131").Value = 2836.5
$ws.Range("N131").Value = -382781.742

# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 1196.0667
$ws.Range("I132").Value = 1099.8334
$ws.Range("K132").Value = 9898.500599999999
$ws.Range("M132").Value = -7368.500599999999

# Row 138: Bring Me Your Tacos / Tacos Al Pastor
$ws.Range("H138").Value = 95110.16
$ws.Range("I138").Value = 1488.3334
$ws.Range("J138").Value = 151283.25
$ws.Range("K138").Value = 4465.0002
$ws.Range("L138").Value = 453849.75
$ws.Range("M138").Value = 674.9997999999996
$ws.Range("N138").Value = -464129.75

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Range("H97").Value = 3138.6
$ws.Range("I97").Value = 1399.3334
$ws.Range("J97").Value = 5747.5
$ws.Range("K97").Value = 1399.3334
$ws.Range("L97").Value = 5747.5
$ws.Range("M97").Value = -903.3334
$ws.Range("N97").Value = -6739.5

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 959.8
$ws.Range("I102").Value = 952.1429000000001
$ws.Range("K102").Value = 952.1429000000001
$ws.Range("M102").Value = 669.8570999999999

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 2748.3333
$ws.Range("I122").Value = 2200
$ws.Range("J122").Value = 3022.5
$ws.Range("K122").Value = 6600
$ws.Range("L122").Value = 9067.5
$ws.Range("M122").Value = -4150
$ws.Range("N122").Value = -13967.5

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 4229.032
$ws.Range("I126").Value = 3145
$ws.Range("J126").Value = 6200
$ws.Range("K126").Value = 9435
$ws.Range("L126").Value = 18600
$ws.Range("M126").Value = -6965
$ws.Range("N126").Value = -23540

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 4260.467
$ws.Range("I7").Value = 4092.8462
$ws.Range("K7").Value = 4092.8462
$ws.Range("M7").Value = -3980.8462

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 5519.5
$ws.Range("I61").Value = 1939.2
$ws.Range("J61").Value = 9099.799999999999
$ws.Range("K61").Value = 1939.2
$ws.Range("L61").Value = 9099.799999999999
$ws.Range("M61").Value = -1737.2
$ws.Range("N61").Value = -9503.799999999999

# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 4036.6316
$ws.Range("I68").Value = 2269.7
$ws.Range("J68").Value = 5999.8887
$ws.Range("K68").Value = 2269.7
$ws.Range("L68").Value = 5999.8887
$ws.Range("M68").Value = -1520.7
$ws.Range("N68").Value = -7497.8887

# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 4036.6316
$ws.Range("I71").Value = 2269.7
$ws.Range("J71").Value = 5999.8887
$ws.Range("K71").Value = 11348.5
$ws.Range("L71").Value = 29999.4435
$ws.Range("M71").Value = -7604.5
$ws.Range("N71").Value = -37487.4435

# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 2898.5715
$ws.Range("I93").Value = 2898.5715
$ws.Range("K93").Value = 2898.5715
$ws.Range("M93").Value = -1650.5715

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 5519.5
$ws.Range("I113").Value = 1939.2
$ws.Range("J113").Value = 9099.799999999999
$ws.Range("K113").Value = 1939.2
$ws.Range("L113").Value = 9099.799999999999
$ws.Range("M113").Value = 230.8
$ws.Range("N113").Value = -13439.8

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 1964315.2
$ws.Range("I122").Value = 2453681
$ws.Range("J122").Value = 6852.5
$ws.Range("K122").Value = 7361043
$ws.Range("L122").Value = 20557.5
$ws.Range("M122").Value = -7358593
$ws.Range("N122").Value = -25457.5

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 4260.467
$ws.Range("I126").Value = 4092.8462
$ws.Range("K126").Value = 12278.5386
$ws.Range("M126").Value = -9808.5386

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 2258.3
$ws.Range("I132").Value = 1607.5555
$ws.Range("K132").Value = 4822.666499999999
$ws.Range("M132").Value = -2292.666499999999

$ws = $wb.Worksheets.Item("WVR")
# Row 74: Clothing the Naked Truth / Ramie Robe of Casting
$ws.Range("H74").Value = 32479
$ws.Range("J74").Value = 32479
$ws.Range("L74").Value = 32479
$ws.Range("N74").Value = -34351

# Row 77: When in Robes (L) / Ramie Robe of Casting
$ws.Range("H77").Value = 32479
$ws.Range("J77").Value = 32479
$ws.Range("L77").Value = 97437
$ws.Range("N77").Value = -106797

# Row 100: Of Great Import / Kudzu Thread
$ws.Range("H100").Value = 1188.3846
$ws.Range("I100").Value = 741.6667
$ws.Range("J100").Value = 1571.2858
$ws.Range("K100").Value = 1483.3334
$ws.Range("L100").Value = 3142.5716
$ws.Range("M100").Value = -942.3334
$ws.Range("N100").Value = -4224.5716

# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 1591.5333
$ws.Range("I122").Value = 1588.2727
$ws.Range("J122").Value = 1600.5
$ws.Range("K122").Value = 4764.8181
$ws.Range("L122").Value = 4801.5
$ws.Range("M122").Value = -2314.8181
$ws.Range("N122").Value = -9701.5

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 37039148
$ws.Range("I136").Value = 55557520
$ws.Range("J136").Value = 2400.4443
$ws.Range("K136").Value = 166672560
$ws.Range("L136").Value = 7201.3329
$ws.Range("M136").Value = -166670010
$ws.Range("N136").Value = -12301.3329
